$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Lookup table: "date|event name" -> new "want-to-go" count (F col)
# ---------------------------------------------------------------
$map = @{}
$map["2024-06-29|苏州·国乙ony茶话会一对一委托-星渡咖啡（取消）"] = 3161
$map["2024-06-30|张家港·幻想物语动漫游戏嘉年华02"] = 1164
$map["2024-06-30|苏州·OrangeOrange夏日随舞派对【免费展会】"] = 141
$map["2024-07-05|苏州·第一届暑假动漫展Comic"] = 127
$map["2024-07-06|苏州·第一届寒假动漫展宅舞比赛-CF01"] = 281
$map["2024-07-12|苏州·暑假动漫展-COS动漫节"] = 100
$map["2024-07-13|张家港·突破次元壁动漫游戏嘉年华"] = 1226
$map["2024-07-19|苏州·萤火国潮文化节动漫品牌博览会"] = 17316
$map["2024-07-20|【大会员提前抢】苏州·OCG国潮动漫游戏嘉年华·你的欲梦内场·全网内场首签"] = 317
$map["2024-07-20|【大会员提前抢】苏州·coke老师撸猫内场票-萤火国潮文化节"] = 222
$map["2024-07-20|张家港·元气爆炸·随机宅舞"] = 1052
$map["2024-07-20|苏州·OCG国潮动漫嘉年华"] = 6601
$map["2024-07-20|苏州·白日梦想7.20全职猎人ONLY展"] = 670
$map["2024-07-21|【大会员提前抢】苏州·OCG国潮动漫游戏嘉年华·火只木南内场票"] = 145
$map["2024-07-21|常熟·SL动漫展02"] = 93
$map["2024-07-26|苏州·暑假COS动漫展-CF01"] = 45
$map["2024-07-27|苏州·AME动漫嘉年华"] = 135
$map["2024-07-27|苏州·第一届动漫游戏展"] = 1287
$map["2024-07-27|苏州·第五人格only·盛典"] = 140
$map["2024-07-27|苏州·音游文化动漫节"] = 47
$map["2024-07-28|张家港·喵言动漫游戏嘉年华"] = 647
$map["2024-07-28|苏州·THSP 05"] = 28
$map["2024-07-28|苏州·第一届维度创想动漫嘉年华"] = 22
$map["2024-08-02|苏州·环球港动漫节"] = 27
$map["2024-08-03|常熟·ACG动漫游戏嘉年华"] = 251
$map["2024-08-03|常熟·CDW.动漫展05"] = 931
$map["2024-08-03|苏州·代号鸢only茶话会-星渡咖啡"] = 90
$map["2024-08-03|苏州·星部落动漫嘉年华"] = 5103
$map["2024-08-04|苏州·授渔仲夏动漫节2.0"] = 522
$map["2024-08-10|苏州·爱乐之城·经典电影作品音乐会"] = 25
$map["2024-08-11|昆山·第七届·xcy新次元动漫嘉年华-狂欢盛典"] = 53
$map["2024-08-17|苏州·ICAN summer World动漫品牌夏游节"] = 11721
$map["2024-08-17|苏州·第二届Redamancy动漫游戏嘉年华"] = 1263
$map["2024-08-24|苏州·赛马娘ONLY"] = 30
$map["2024-09-15|苏州·Good jump ACG中秋嘉年华动漫国潮文化节"] = 180
$map["2024-10-01|苏州·I COME ACG动漫品牌博览会"] = 250
$map["2024-10-01|苏州·第十三届理想乡动漫展-同人创作者大会"] = 3882
$map["2024-10-02|苏州·明日方舟ONLY#2024~佑桑柔"] = 281
$map["2024-10-26|苏州·第三届华盟国漫次元嘉年华"] = 83

# ---------------------------------------------------------------
# Sheet "展览" (exhibitions): remove duplicate rows, then refresh
# the A (index) and F (want-to-go count) columns.
# ---------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

# Duplicate rows to remove, original row numbers, deleted bottom-up
# so earlier indices stay valid while iterating.
$dupRows = @(40, 25, 19, 7, 5, 3)
foreach ($r in $dupRows) {
    $wsExpo.Rows.Item($r).Delete()
}

$lastRowExpo = $wsExpo.Cells.Item(1, 1).Value2
$usedRange = $wsExpo.UsedRange
$lastRowExpo = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRowExpo; $r++) {
    $wsExpo.Cells.Item($r, 1).Value2 = $r - 1
    $date = $wsExpo.Cells.Item($r, 2).Value2
    $name = $wsExpo.Cells.Item($r, 3).Value2
    $key = "$date|$name"
    if ($map.ContainsKey($key)) {
        $wsExpo.Cells.Item($r, 6).Value2 = $map[$key]
    }
}

# ---------------------------------------------------------------
# Sheet "演出" (performances): update want-to-go count for row 2.
# ---------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Cells.Item(2, 6).Value2 = 25

# ---------------------------------------------------------------
# Sheet "全部类型" (all types): refresh F column for every data row,
# and fix the G3 sale-status text.
# ---------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$usedRangeAll = $wsAll.UsedRange
$lastRowAll = $usedRangeAll.Rows.Count

for ($r = 2; $r -le $lastRowAll; $r++) {
    $date = $wsAll.Cells.Item($r, 2).Value2
    $name = $wsAll.Cells.Item($r, 3).Value2
    $key = "$date|$name"
    if ($map.ContainsKey($key)) {
        $wsAll.Cells.Item($r, 6).Value2 = $map[$key]
    }
}

$wsAll.Cells.Item(3, 7).Value2 = "不可售"
